# Generate Report for Handback
#
# The CI "handback" step marks the zh-cn and de-de localization rows as
# handed back: the Status text changes workbook-wide, each language sheet
# gains "Latest Target File" / "Latest Handback File" columns (F/G) with
# hyperlinks mirroring the existing source/handoff links, and the
# "Latest Handback DateTime" column (H) is stamped with the real handback
# time instead of the zero-date placeholder.

function Get-HyperlinkAddress($ws, $addr) {
  foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq $addr) {
      return $hl.Address
    }
  }
  return $null
}

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates (shared text across Overview + both language sheets) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# --- zh-cn sheet: Latest Target File (F) / Latest Handback File (G) / Latest Handback DateTime (H) ---
$zhMdUrl = Get-HyperlinkAddress $wsZh '$A$2'
$zhXlfUrl = Get-HyperlinkAddress $wsZh '$D$2'
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("F2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdUrl, "", "", "a.md")

$wsZh.Range("G2").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)

$wsZh.Range("H2").Value = "2016-03-21 04:31:15"

$wsZh.Range("F3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdUrl, "", "", "a.md")

$wsZh.Range("G3").Value = $zhXlfName
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

$wsZh.Range("H3").Value = "2016-03-21 04:31:15"

# --- de-de sheet: Latest Target File (F) / Latest Handback File (G) / Latest Handback DateTime (H) ---
$deMdUrl = Get-HyperlinkAddress $wsDe '$A$2'
$deXlfUrl = Get-HyperlinkAddress $wsDe '$D$2'
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("F2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdUrl, "", "", "a.md")

$wsDe.Range("G2").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)

$wsDe.Range("H2").Value = "2016-03-21 04:31:21"

$wsDe.Range("F3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdUrl, "", "", "a.md")

$wsDe.Range("G3").Value = $deXlfName
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$wsDe.Range("H3").Value = "2016-03-21 04:31:21"

Write-Host "Handback report generated."
